$d = $word.ActiveDocument

function Get-ParagraphIndexContaining($doc, $searchText) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text.Contains($searchText)) {
            return $i
        }
    }
    return -1
}

# -----------------------------------------------------------------------
# Edit 1: merge the " " run and the "What was their..." run that follow
# "Why did the trialist perform a non-inferiority trial instead of a
# superiority trial?" into a single run (no visible text change).
# -----------------------------------------------------------------------
$idx1 = Get-ParagraphIndexContaining $d "Why did the trialist perform"
$p1 = $d.Paragraphs($idx1)
$r1 = $p1.Range
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:u w:val="single"/></w:rPr><w:t>Why did the trialist perform a non-inferiority trial instead of a superiority trial?</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> What was their non-inferiority margin, and is it justified? If the intervention is found non-inferior to the control, what advantages would result? Does this benefit outweigh the possible inferiority margin?</w:t></w:r></w:p>'
$r1.InsertXML($xml1)

# -----------------------------------------------------------------------
# Edit 2: remove the "_GoBack" bookmark sitting between the two runs in
# the "What is a noninferiority margin?" paragraph, and merge those two
# runs into one (no visible text change).
# -----------------------------------------------------------------------
$idx2 = Get-ParagraphIndexContaining $d "What is a noninferiority margin?"
$p2 = $d.Paragraphs($idx2)
$r2 = $p2.Range
$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:u w:val="single"/></w:rPr><w:t>What is a noninferiority margin?</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve"> The </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>amo</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>unt</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve"> less</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve"> that the trial investigator is willing to tolerate and still call the treatments equal. This should have some relation to what the minimally important clinical difference is. Importantly, this can’t be 0 (or else, we’d have shown superiority). However, it should be closer to 0 if the outcome is severe (meaning, we’d tolerate very little extra mortality) or if the benefits of the comparator treatment are not that great (meaning, if a drug that doesn’t have many other benefits much inferior, that’s not helpful) </w:t></w:r></w:p>'
$r2.InsertXML($xml2)

# -----------------------------------------------------------------------
# Edit 3: append, at the very end of the document:
#   - an empty paragraph
#   - a paragraph with "Good layman targeted summary: " followed by a
#     hyperlink to the nephjc.com article
#   - an empty paragraph containing the (moved) "_GoBack" bookmark
# -----------------------------------------------------------------------
$end = $d.Content
$end.Collapse(0)
$appendXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:rPr><w:rFonts w:ascii='Calibri' w:hAnsi='Calibri' w:cs='Calibri'/></w:rPr></w:pPr></w:p>" + `
             "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:rPr><w:rFonts w:ascii='Calibri' w:hAnsi='Calibri' w:cs='Calibri'/></w:rPr><w:t xml:space='preserve'>Good layman targeted summary: </w:t></w:r></w:p>" + `
             "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'></w:p>" + `
             "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:rPr><w:rFonts w:ascii='Calibri' w:hAnsi='Calibri' w:cs='Calibri'/></w:rPr></w:pPr></w:p>"
$end.InsertXML($appendXml)

# Add the hyperlink into the (currently empty) placeholder paragraph.
$hlParaIndex = $d.Paragraphs.Count - 1
$hlPara = $d.Paragraphs($hlParaIndex)
[void]$d.Hyperlinks.Add($hlPara.Range, "http://www.nephjc.com/news/2019/7/8/understanding-the-vortex-of-non-inferiority-trials", "", "", "http://www.nephjc.com/news/2019/7/8/understanding-the-vortex-of-non-inferiority-trials")

# Merge the "Good layman..." paragraph with the hyperlink paragraph by
# deleting the paragraph mark between them.
$glParaIndex = $hlParaIndex - 1
$glPara = $d.Paragraphs($glParaIndex)
$joinRange = $glPara.Range
$joinRange.Collapse(0)
$joinRange.MoveEnd(1, 1)
$joinRange.Delete()

# Put the "_GoBack" bookmark in the final (last) empty paragraph.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
[void]$d.Bookmarks.Add("_GoBack", $lastPara.Range)
